$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-18 (Item name, Price, Shop)
$data = @(
    @("GeForce RTX 3060 Ventus 2X OC 12GB GDDR6 192bit MSI", 290, "RD Electronics"),
    @("GeForce RTX 3060 Gaming OC 8GB GDDR6 128bit rev. 2.0 Gigabyte", 292, "RD Electronics"),
    @("GeForce RTX 3060 Gaming OC 8GB GDDR6 128bit Gigabyte", 292.71, "RD Electronics"),
    @("GeForce RTX 3060 WindForce OC 12GB GDDR6 192bit Gigabyte", 300.2, "RD Electronics"),
    @("GeForce RTX 3060 Gaming OC 12GB GDDR6 192bit rev. 2.0 LHR Gigabyte", 302.27, "RD Electronics"),
    @("GeForce RTX 3060 Dual OC V2 LHR 12GB GDDR6 192bit Asus", 317.9, "RD Electronics"),
    @("GeForce RTX Extreme Gaming 3060 12GB GDDR6 192bit Biostar", 317.9, "RD Electronics"),
    @("GeForce RTX 3060 V2 OC 12GB GDDR6 192bit Asus", 318.08, "RD Electronics"),
    @("GeForce RTX 3060 Twin Twin Edge 12GB GDDR6 192bit ZOTAC", 334.98, "RD Electronics"),
    @("GeForce RTX 3060 Twin X2 OC LHR 12GB GDDR6 192bit Inno3D", 341, "RD Electronics"),
    @("GeForce RTX 3060 Twin X2 12GB GDDR6 192bit Inno3D", 351, "RD Electronics"),
    @("GeForce RTX 3060 Ti Extreme Gaming 8GB GDDR6X 256bit Biostar", 361, "RD Electronics"),
    @("Ryzen 5 5500 16GB 250SSD RTX3060 12GB no-OS INTOP", 705.99, "RD Electronics"),
    @("Ryzen 5 5500 16GB 500SSD RTX3060 12GB no-OS INTOP", 718.99, "RD Electronics"),
    @("Ryzen 5 5600X 16GB 250SSD RTX3060 12GB no-OS INTOP", 740.99, "RD Electronics"),
    @("Ryzen 5 5500 16GB 1TB SSD RTX3060 12GB no-OS INTOP", 740.99, "RD Electronics"),
    @("Blade RZ09-0421EED3-R3E1 15.6 240hz i7-12800H 16GB 1SSD RTX3060 EN W11 Razer", 3340, "RD Electronics")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row = $row + 1
}

# Remove the now-unused rows 19:31 (sheet shrinks from 31 to 18 rows)
$ws.Rows("19:31").Delete()
